# Update "想去人数" (interest count) figures in column F on the
# "展览" (Exhibition) and "全部类型" (All types) sheets, as the site data
# was refreshed at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" -> row number : new value
$exhibitionUpdates = @{
    14 = 434
    15 = 1357
    17 = 108
    20 = 36
    22 = 1004
    25 = 22
    26 = 5878
    29 = 99
    31 = 14531
    32 = 1439
    33 = 212
    36 = 9068
    37 = 625
    39 = 144
}

# Sheet "全部类型" -> row number : new value
$allTypesUpdates = @{
    14 = 434
    15 = 1357
    17 = 108
    21 = 36
    24 = 1004
    27 = 22
    29 = 5878
    32 = 99
    34 = 14532
    35 = 1439
    36 = 212
    39 = 9068
    40 = 625
    42 = 144
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
